$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05165006512414756
$ws.Range("B3").Value = -0.21663876924289746
$ws.Range("B4").Value = 2.0352586040178022

$ws.Rows.Item(5).Delete()
